$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting the existing row 9 (now row 10) down
$ws.Rows.Item(9).Insert()

# New row 9: "Niet te lokaliseren"
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 99993
$ws.Range("C9").Value = "Niet te lokaliseren"
$ws.Range("D9").Value = "Niet te lokaliseren"

# Row 10 (previously row 9): update volgnr from 8 to 9
$ws.Range("A10").Value = 9
